$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.984.21'
$ws.Range("E2").Value = '  -1.28%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.497.61'
$ws.Range("E3").Value = '  -0.19%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '581.67'
$ws.Range("E5").Value = '  -0.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '133.19'
$ws.Range("E6").Value = '  -0.29%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.496.81'
$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.488'
$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.123'
$ws.Range("E10").Value = '  -0.80%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.05'
$ws.Range("E11").Value = '  -1.12%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.372'
$ws.Range("E12").Value = '  -2.65%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.092.91'
$ws.Range("E13").Value = '  -0.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.15'
$ws.Range("E14").Value = '  -1.10%  '

$ws.Range("E15").Value = '  +1.40%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.498.74'
$ws.Range("E16").Value = '  -0.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000177'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.045.03'
$ws.Range("E18").Value = '  -1.21%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.64'
$ws.Range("E19").Value = '  -3.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.77'
$ws.Range("E20").Value = '  -2.92%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.55'
$ws.Range("E21").Value = '  -1.09%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '380.60'
$ws.Range("E22").Value = '  -2.21%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.640.84'
$ws.Range("E23").Value = '  -0.23%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.563'
$ws.Range("E24").Value = '  -1.33%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '73.94'
$ws.Range("E25").Value = '  -0.58%  '

$ws.Range("E26").Value = '  -0.01%  '

$ws.Range("E27").Value = '  +0.57%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000113'
$ws.Range("E28").Value = '  +3.25%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.55'
$ws.Range("E29").Value = '  -1.70%  '

$ws.Range("E30").Value = '  +0.43%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.27'
$ws.Range("E32").Value = '  +1.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.20'
$ws.Range("E33").Value = '  -2.00%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.513.04'
$ws.Range("E34").Value = '  +0.10%  '

$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '23.42'
$ws.Range("E36").Value = '  -2.01%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.144'
$ws.Range("E37").Value = '  -0.51%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.29'
$ws.Range("E38").Value = '  +3.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.86'
$ws.Range("E39").Value = '  -0.47%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.53'
$ws.Range("E40").Value = '  -1.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '162.02'
$ws.Range("E41").Value = '  -3.88%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0779'
$ws.Range("E42").Value = '  -2.24%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '26.39'
$ws.Range("E43").Value = '  +1.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.804'
$ws.Range("E44").Value = '  -1.18%  '

$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '41.53'
$ws.Range("E46").Value = '  -2.32%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.19'
$ws.Range("E47").Value = '  -3.93%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.35'
$ws.Range("E48").Value = '  -0.22%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.57'
$ws.Range("E49").Value = '  -4.21%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.460.98'
$ws.Range("E50").Value = '  -0.03%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.73'
$ws.Range("E51").Value = '  -1.21%  '
